$wb = $excel.ActiveWorkbook

# Update the "Data" sheet: row 5's browser (column C) changes from "firefox" to "chrome"
$ws = $wb.Worksheets.Item("Data")
$ws.Range("C5").Value = "chrome"

# Move the active selection to D5 on the Data sheet (matches the cursor state in the target file)
$ws.Activate()
$ws.Range("D5").Select()
